# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, shared between the two sheets that
# carry the exhibition ("展览") data: the main "展览" sheet and the combined
# "全部类型" sheet. Note row 35 in "展览" corresponds to row 37 in
# "全部类型" (the "全部类型" sheet has one extra preceding row for this
# particular event), and row 39 in "展览" corresponds to row 41 in
# "全部类型".

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 15703
$wsExhibition.Range("F8").Value = 707
$wsExhibition.Range("F9").Value = 15440
$wsExhibition.Range("F10").Value = 56
$wsExhibition.Range("F11").Value = 9036
$wsExhibition.Range("F15").Value = 95
$wsExhibition.Range("F18").Value = 200
$wsExhibition.Range("F21").Value = 551
$wsExhibition.Range("F26").Value = 3
$wsExhibition.Range("F29").Value = 89
$wsExhibition.Range("F35").Value = 255
$wsExhibition.Range("F39").Value = 5561

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F5").Value = 15703
$wsAllTypes.Range("F8").Value = 707
$wsAllTypes.Range("F9").Value = 15440
$wsAllTypes.Range("F10").Value = 56
$wsAllTypes.Range("F11").Value = 9036
$wsAllTypes.Range("F15").Value = 95
$wsAllTypes.Range("F18").Value = 200
$wsAllTypes.Range("F21").Value = 551
$wsAllTypes.Range("F26").Value = 3
$wsAllTypes.Range("F29").Value = 89
$wsAllTypes.Range("F37").Value = 255
$wsAllTypes.Range("F41").Value = 5561
